$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 8 more match rows (10-17) to the "Kagiso Rabada" log sheet.
# These are re-listings of the existing 8 matches (rows 2-9) in a new
# order, mirroring the upstream JSON -> xlsx re-export.
# Columns G-K look numeric but are stored as text in the source file,
# so they are written here with a leading apostrophe to force Excel to
# keep them as literal text instead of coercing to numbers.

# row 10 (duplicate of source row 8)
$ws.Cells.Item(10,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(10,2).Value = ' November 05 2020'
$ws.Cells.Item(10,3).Value = 'Mumbai won by 57 runs'
$ws.Cells.Item(10,4).Value = 'Delhi Capitals'
$ws.Cells.Item(10,5).Value = 'Mumbai Indians'
$ws.Cells.Item(10,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(10,7).Value = '''15'
$ws.Cells.Item(10,8).Value = '''15'
$ws.Cells.Item(10,9).Value = '''2'
$ws.Cells.Item(10,10).Value = '''0'
$ws.Cells.Item(10,11).Value = '''100.00'

# row 11 (duplicate of source row 9)
$ws.Cells.Item(11,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(11,2).Value = ' October 27 2020'
$ws.Cells.Item(11,3).Value = 'Sunrisers won by 88 runs'
$ws.Cells.Item(11,4).Value = 'Delhi Capitals'
$ws.Cells.Item(11,5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(11,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(11,7).Value = '''3'
$ws.Cells.Item(11,8).Value = '''7'
$ws.Cells.Item(11,9).Value = '''0'
$ws.Cells.Item(11,10).Value = '''0'
$ws.Cells.Item(11,11).Value = '''42.85'

# row 12 (duplicate of source row 4)
$ws.Cells.Item(12,1).Value = ' Abu Dhabi'
$ws.Cells.Item(12,2).Value = ' September 29 2020'
$ws.Cells.Item(12,3).Value = 'Sunrisers won by 15 runs'
$ws.Cells.Item(12,4).Value = 'Delhi Capitals'
$ws.Cells.Item(12,5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(12,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(12,7).Value = '''15'
$ws.Cells.Item(12,8).Value = '''7'
$ws.Cells.Item(12,9).Value = '''1'
$ws.Cells.Item(12,10).Value = '''1'
$ws.Cells.Item(12,11).Value = '''214.28'

# row 13 (duplicate of source row 7)
$ws.Cells.Item(13,1).Value = ' Abu Dhabi'
$ws.Cells.Item(13,2).Value = ' October 24 2020'
$ws.Cells.Item(13,3).Value = 'KKR won by 59 runs'
$ws.Cells.Item(13,4).Value = 'Delhi Capitals'
$ws.Cells.Item(13,5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(13,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(13,7).Value = '''9'
$ws.Cells.Item(13,8).Value = '''10'
$ws.Cells.Item(13,9).Value = '''1'
$ws.Cells.Item(13,10).Value = '''0'
$ws.Cells.Item(13,11).Value = '''90.00'

# row 14 (duplicate of source row 5)
$ws.Cells.Item(14,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(14,2).Value = ' October 31 2020'
$ws.Cells.Item(14,3).Value = 'Mumbai won by 9 wickets (with 34 balls remaining)'
$ws.Cells.Item(14,4).Value = 'Delhi Capitals'
$ws.Cells.Item(14,5).Value = 'Mumbai Indians'
$ws.Cells.Item(14,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(14,7).Value = '''12'
$ws.Cells.Item(14,8).Value = '''7'
$ws.Cells.Item(14,9).Value = '''0'
$ws.Cells.Item(14,10).Value = '''1'
$ws.Cells.Item(14,11).Value = '''171.42'

# row 15 (duplicate of source row 6)
$ws.Cells.Item(15,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(15,2).Value = ' November 10 2020'
$ws.Cells.Item(15,3).Value = 'Mumbai won by 5 wickets (with 8 balls remaining)'
$ws.Cells.Item(15,4).Value = 'Delhi Capitals'
$ws.Cells.Item(15,5).Value = 'Mumbai Indians'
$ws.Cells.Item(15,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(15,7).Value = '''0'
$ws.Cells.Item(15,8).Value = '''0'
$ws.Cells.Item(15,9).Value = '''0'
$ws.Cells.Item(15,10).Value = '''0'
$ws.Cells.Item(15,11).Value = '''-'

# row 16 (duplicate of source row 2)
$ws.Cells.Item(16,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(16,2).Value = ' September 20 2020'
$ws.Cells.Item(16,3).Value = 'Match tied (Capitals won the one-over eliminator)'
$ws.Cells.Item(16,4).Value = 'Delhi Capitals'
$ws.Cells.Item(16,5).Value = 'Kings XI Punjab'
$ws.Cells.Item(16,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(16,7).Value = '''0'
$ws.Cells.Item(16,8).Value = '''0'
$ws.Cells.Item(16,9).Value = '''0'
$ws.Cells.Item(16,10).Value = '''0'
$ws.Cells.Item(16,11).Value = '''-'

# row 17 (duplicate of source row 3)
$ws.Cells.Item(17,1).Value = ' Sharjah'
$ws.Cells.Item(17,2).Value = ' October 09 2020'
$ws.Cells.Item(17,3).Value = 'Capitals won by 46 runs'
$ws.Cells.Item(17,4).Value = 'Delhi Capitals'
$ws.Cells.Item(17,5).Value = 'Rajasthan Royals'
$ws.Cells.Item(17,6).Value = 'Kagiso Rabada '
$ws.Cells.Item(17,7).Value = '''2'
$ws.Cells.Item(17,8).Value = '''3'
$ws.Cells.Item(17,9).Value = '''0'
$ws.Cells.Item(17,10).Value = '''0'
$ws.Cells.Item(17,11).Value = '''66.66'

